$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 543, pushing all the
# existing records from row 543 down to row 544 (through the former
# last row 594, which becomes row 595). Replicate that with a native
# row insert, then populate the freshly inserted row with the new data.
$ws.Rows.Item(543).Insert()

$ws.Range("A543").Value = 3
$ws.Range("B543").Value = "Femacal de La Calera"
$ws.Range("C543").Value = "Coquimbo"
$ws.Range("D543").Value = 45132
$ws.Range("E543").Value = 5
$ws.Range("F543").Value = 100114013
$ws.Range("G543").Value = "Zanahoria"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 160
$ws.Range("K543").Value = 7500
$ws.Range("L543").Value = 7500
$ws.Range("M543").Value = 7500
$ws.Range("N543").Value = "`$/saco 20 kilos"
$ws.Range("O543").Value = "Provincia de Quillota"
$ws.Range("P543").Value = 375
$ws.Range("Q543").Value = 20
$ws.Range("R543").Value = "Hortaliza"
